# Updates pricing/profit figures across several Leve sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) to reflect refreshed market board averages.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 2069.6
$ws.Range("J32").Value = 2069.6
$ws.Range("H32").Value = 1974.6666
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -2721.6
$ws.Range("I32").Value = 1500

$ws.Range("L80").Value = 780.5999999999999
$ws.Range("H80").Value = 1839.5454
$ws.Range("N80").Value = -2776.6
$ws.Range("J80").Value = 260.2

$ws.Range("L83").Value = 2341.8
$ws.Range("H83").Value = 1839.5454
$ws.Range("J83").Value = 260.2
$ws.Range("N83").Value = -12325.8

$ws.Range("I86").Value = 1579.4
$ws.Range("K86").Value = 1579.4
$ws.Range("M86").Value = -456.4000000000001
$ws.Range("H86").Value = 1617.6364

$ws.Range("H89").Value = 1617.6364
$ws.Range("M89").Value = -2281
$ws.Range("K89").Value = 7897
$ws.Range("I89").Value = 1579.4

$ws.Range("J138").Value = 3216.106
$ws.Range("L138").Value = 9648.318000000001
$ws.Range("N138").Value = -19928.318
$ws.Range("H138").Value = 3490.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 14690.467
$ws.Range("H32").Value = 18266.795
$ws.Range("M32").Value = -14403.467
$ws.Range("I32").Value = 14690.467

$ws.Range("K61").Value = 33283.4
$ws.Range("I61").Value = 33283.4
$ws.Range("H61").Value = 26918.828
$ws.Range("M61").Value = -33071.4

$ws.Range("K63").Value = 1133
$ws.Range("M63").Value = -447
$ws.Range("H63").Value = 2569.077
$ws.Range("I63").Value = 1133

$ws.Range("M66").Value = -2233
$ws.Range("K66").Value = 5665
$ws.Range("H66").Value = 2569.077
$ws.Range("I66").Value = 1133

$ws.Range("I97").Value = 1158.1428
$ws.Range("M97").Value = -662.1428000000001
$ws.Range("H97").Value = 1381.5
$ws.Range("K97").Value = 1158.1428

$ws.Range("I102").Value = 1750
$ws.Range("K102").Value = 1750
$ws.Range("M102").Value = -128
$ws.Range("H102").Value = 1750

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("K136").Value = 99850.20000000001
$ws.Range("I136").Value = 33283.4
$ws.Range("H136").Value = 26918.828
$ws.Range("M136").Value = -97300.20000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J20").Value = 4099.3335
$ws.Range("L20").Value = 4099.3335
$ws.Range("N20").Value = -4593.3335
$ws.Range("H20").Value = 2766.4211

$ws.Range("J22").Value = 997.5
$ws.Range("L22").Value = 997.5
$ws.Range("H22").Value = 773.25
$ws.Range("N22").Value = -1343.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4306.2593
$ws.Range("I31").Value = 2269
$ws.Range("K31").Value = 2269
$ws.Range("M31").Value = -1974

$ws.Range("I34").Value = 2269
$ws.Range("K34").Value = 2269
$ws.Range("M34").Value = -2067
$ws.Range("H34").Value = 4306.2593

$ws.Range("J122").Value = 1264
$ws.Range("N122").Value = -8692
$ws.Range("L122").Value = 3792
$ws.Range("H122").Value = 1386

$ws.Range("J132").Value = 5137
$ws.Range("M132").Value = -870.7999999999997
$ws.Range("I132").Value = 1133.6
$ws.Range("L132").Value = 15411
$ws.Range("K132").Value = 3400.8
$ws.Range("N132").Value = -20471
$ws.Range("H132").Value = 1800.8334

$ws.Range("H134").Value = 1026.5686
$ws.Range("K134").Value = 2583.5853
$ws.Range("M134").Value = -48.58530000000019
$ws.Range("I134").Value = 861.1951

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J122").Value = 1520
$ws.Range("M122").Value = -4750
$ws.Range("N122").Value = -18580
$ws.Range("L122").Value = 13680
$ws.Range("H122").Value = 1389.091
$ws.Range("K122").Value = 7200
$ws.Range("I122").Value = 800

$ws.Range("J140").Value = 3575.7058
$ws.Range("K140").Value = 2784.9231
$ws.Range("H140").Value = 2428.5
$ws.Range("N140").Value = -21087.1174
$ws.Range("L140").Value = 10727.1174
$ws.Range("M140").Value = 2395.0769
$ws.Range("I140").Value = 928.3077

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 49999
$ws.Range("N15").Value = -50575
$ws.Range("J15").Value = 49999
$ws.Range("L15").Value = 49999

$ws.Range("L46").Value = 28250
$ws.Range("N46").Value = -28562
$ws.Range("J46").Value = 28250
$ws.Range("H46").Value = 28250

$ws.Range("L80").Value = 1378
$ws.Range("H80").Value = 2356.0588
$ws.Range("N80").Value = -3374
$ws.Range("I80").Value = 2657
$ws.Range("J80").Value = 1378
$ws.Range("M80").Value = -1659
$ws.Range("K80").Value = 2657

$ws.Range("J81").Value = 49999
$ws.Range("N81").Value = -51995
$ws.Range("L81").Value = 49999
$ws.Range("H81").Value = 49999

$ws.Range("L83").Value = 6890
$ws.Range("H83").Value = 2356.0588
$ws.Range("J83").Value = 1378
$ws.Range("N83").Value = -16874
$ws.Range("M83").Value = -8293
$ws.Range("K83").Value = 13285
$ws.Range("I83").Value = 2657

$ws.Range("N84").Value = -159981
$ws.Range("J84").Value = 49999
$ws.Range("H84").Value = 49999
$ws.Range("L84").Value = 149997

$ws.Range("I102").Value = 2063
$ws.Range("L102").Value = 2940
$ws.Range("K102").Value = 2063
$ws.Range("J102").Value = 2940
$ws.Range("N102").Value = -6184
$ws.Range("M102").Value = -441
$ws.Range("H102").Value = 2306.611

$ws.Range("M122").Value = -1948.6
$ws.Range("H122").Value = 1591.25
$ws.Range("K122").Value = 4398.6
$ws.Range("I122").Value = 1466.2

$ws.Range("J132").Value = 2627
$ws.Range("M132").Value = -3215840.3
$ws.Range("I132").Value = 1072790.1
$ws.Range("L132").Value = 7881
$ws.Range("K132").Value = 3218370.3
$ws.Range("N132").Value = -12941
$ws.Range("H132").Value = 678519.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -2598.1
$ws.Range("K7").Value = 2710.1
$ws.Range("I7").Value = 2710.1
$ws.Range("H7").Value = 3005.2173

$ws.Range("L23").Value = 0
$ws.Range("H23").Value = 10000000
$ws.Range("J23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("N43").Value = -62788.8
$ws.Range("H43").Value = 62402.8
$ws.Range("L43").Value = 62402.8
$ws.Range("J43").Value = 62402.8

$ws.Range("J109").Value = 59996
$ws.Range("N109").Value = -62770
$ws.Range("L109").Value = 59996
$ws.Range("H109").Value = 59996

$ws.Range("K126").Value = 8130.299999999999
$ws.Range("H126").Value = 3005.2173
$ws.Range("I126").Value = 2710.1
$ws.Range("M126").Value = -5660.299999999999

$ws.Range("J132").Value = 4898.5
$ws.Range("M132").Value = -7714.0355
$ws.Range("I132").Value = 3414.6785
$ws.Range("L132").Value = 14695.5
$ws.Range("K132").Value = 10244.0355
$ws.Range("N132").Value = -19755.5
$ws.Range("H132").Value = 4206.05

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J123").Value = 49053.4
$ws.Range("N123").Value = -58853.4
$ws.Range("H123").Value = 49053.4
$ws.Range("L123").Value = 49053.4
